$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert 6 new rows right before the old row 68 (河北卫视 block),
# shifting everything else down by 6 rows. The new channel (海峡卫视)
# sorts alphabetically between 海南卫视 and 河北卫视.
$ws.Range("A68:A73").EntireRow.Insert()

$ipLabels = @("IPA", "IPB", "IPC", "IPD", "IPE", "IPF")
for ($i = 0; $i -lt 6; $i++) {
    $r = 68 + $i
    $ws.Cells.Item($r, 1).Value = "海峡卫视"
    $ws.Cells.Item($r, 2).Value = ","
    $ws.Cells.Item($r, 3).Value = "http://"
    $ws.Cells.Item($r, 4).Value = $ipLabels[$i]
    $ws.Cells.Item($r, 5).Value = "/bst/hxwshd4m@4000000.m3u8"
}

# Refresh the sheet's AutoFilter so its range grows to cover the new rows.
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:F223").AutoFilter()

# The workbook-level hidden _xlnm._FilterDatabase name for this sheet still
# points at the old range; update it to match.
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.RefersTo -like "*F`$217*") {
        $n.RefersTo = "=卫视!`$A`$1:`$F`$223"
    }
}
